# Fruta / hortaliza, semanal
# Re-shuffle the per-row data block (Fecha, Variedad, Calidad, Volumen,
# Precio minimo, Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg)
# across rows 2-28 according to the mapping derived from the target diff.
# Columns A,B,C,E,F,G,N,Q,R stay put; only D,H,I,J,K,L,M,O,P move as a block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (i.e. destinationRow will receive the block
# that currently lives in sourceRow)
$map = @{
    2  = 18
    3  = 14
    4  = 4
    5  = 6
    6  = 10
    7  = 7
    8  = 15
    9  = 25
    10 = 24
    11 = 19
    12 = 12
    13 = 13
    14 = 23
    15 = 22
    16 = 11
    17 = 26
    18 = 8
    19 = 5
    20 = 16
    21 = 17
    22 = 21
    23 = 28
    24 = 9
    25 = 3
    26 = 2
    27 = 20
    28 = 27
}

$cols = @(4, 8, 9, 10, 11, 12, 13, 15, 16)  # D, H, I, J, K, L, M, O, P

# 1) Snapshot every source row's block values BEFORE any writes, so that
#    permutation cycles don't clobber data we still need to read.
$snapshot = @{}
for ($row = 2; $row -le 28; $row++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowVals
}

# 2) Write each destination row's block from the snapshot of its mapped
#    source row.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value = $srcVals[$col]
    }
}
